$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: assign a text-looking numeric value to a cell while keeping it
# stored as text (matching the source inlineStr cells) and without leaving
# the cell on a non-default style.
function Set-TextValue($cellRef, $text) {
    $cell = $ws.Range($cellRef)
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = "Normal"
}

Set-TextValue 'D2' '70.755.74'
$ws.Range('E2').Value = '  -0.28%  '

Set-TextValue 'D3' '3.528.29'
$ws.Range('E3').Value = '  -0.99%  '

Set-TextValue 'D4' '0.999'

Set-TextValue 'D5' '612.44'
$ws.Range('E5').Value = '  -0.58%  '

Set-TextValue 'D6' '173.73'
$ws.Range('E6').Value = '  +0.85%  '

Set-TextValue 'D7' '3.525.10'
$ws.Range('E7').Value = '  -0.98%  '

$ws.Range('E8').Value = '  -1.34%  '

$ws.Range('E9').Value = '  -0.04%  '

Set-TextValue 'D10' '0.197'
$ws.Range('E10').Value = '  -0.52%  '

$ws.Range('E11').Value = '  +1.46%  '

$ws.Range('E12').Value = '  +0.05%  '

Set-TextValue 'D13' '46.62'
$ws.Range('E13').Value = '  -0.27%  '

$ws.Range('E14').Value = '  -0.73%  '

Set-TextValue 'D15' '4.098.72'
$ws.Range('E15').Value = '  -0.85%  '

Set-TextValue 'D16' '8.45'
$ws.Range('E16').Value = '  +0.29%  '

Set-TextValue 'D17' '615.37'
$ws.Range('E17').Value = '  -0.93%  '

Set-TextValue 'D18' '3.528.61'
$ws.Range('E18').Value = '  -0.92%  '

Set-TextValue 'D19' '70.756.91'
$ws.Range('E19').Value = '  -0.29%  '

$ws.Range('E20').Value = '  +1.79%  '

Set-TextValue 'D21' '17.78'
$ws.Range('E21').Value = '  +2.28%  '

Set-TextValue 'D22' '0.886'
$ws.Range('E22').Value = '  +0.32%  '

Set-TextValue 'D23' '8.99'
$ws.Range('E23').Value = '  -5.15%  '

Set-TextValue 'D24' '15.74'
$ws.Range('E24').Value = '  -0.04%  '

Set-TextValue 'D25' '98.18'
$ws.Range('E25').Value = '  +1.26%  '

$ws.Range('E26').Value = '  -1.59%  '

$ws.Range('E27').Value = '  +0.03%  '

Set-TextValue 'D28' '2.61'
$ws.Range('E28').Value = '  -0.40%  '

Set-TextValue 'D29' '33.84'
$ws.Range('E29').Value = '  +0.79%  '

Set-TextValue 'D30' '9.15'
$ws.Range('E30').Value = '  +0.73%  '

Set-TextValue 'D31' '3.04'
$ws.Range('E31').Value = '  -1.28%  '

Set-TextValue 'D32' '8.17'
$ws.Range('E32').Value = '  -4.40%  '

$ws.Range('E33').Value = '  -0.41%  '

Set-TextValue 'D34' '6.86'
$ws.Range('E34').Value = '  -1.76%  '

Set-TextValue 'D35' '615.36'

Set-TextValue 'D36' '0.101'
$ws.Range('E36').Value = '  -0.81%  '

Set-TextValue 'D37' '10.86'
$ws.Range('E37').Value = '  -0.27%  '

Set-TextValue 'D38' '3.53'
$ws.Range('E38').Value = '  -3.09%  '

Set-TextValue 'D39' '0.0476'
$ws.Range('E39').Value = '  -0.18%  '

Set-TextValue 'D40' '57.08'
$ws.Range('E40').Value = '  -1.14%  '

$ws.Range('E41').Value = '  +0.13%  '

Set-TextValue 'D42' '0.144'
$ws.Range('E42').Value = '  +0.80%  '

Set-TextValue 'D43' '3.380.75'
$ws.Range('E43').Value = '  +0.20%  '

Set-TextValue 'D44' '0.0₃0741'

Set-TextValue 'D45' '0.314'
$ws.Range('E45').Value = '  -2.28%  '

$ws.Range('E46').Value = '  -1.89%  '

Set-TextValue 'D47' '32.30'
$ws.Range('E47').Value = '  -2.29%  '

$ws.Range('E49').Value = '  +0.07%  '

Set-TextValue 'D50' '133.98'
$ws.Range('E50').Value = '  +0.13%  '
